# Fruta / hortaliza, semanal
# Insert a new weekly record as row 213 in the Mango price sheet, shifting
# the existing rows 213:265 down to 214:266.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row at position 213 (existing rows 213-265 shift down
# to 214-266, preserving their data/formatting).
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new weekly entry.
$ws.Range("A213").Value = 9
$ws.Range("B213").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C213").Value = "Metropolitana"
$ws.Range("D213").Value = 44476
$ws.Range("D213").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E213").Value = 13
$ws.Range("F213").Value = "Fruta"
$ws.Range("G213").Value = 100108
$ws.Range("H213").Value = "Tropicales y subtropicales"
$ws.Range("I213").Value = 100108002
$ws.Range("J213").Value = "Mango"
$ws.Range("K213").Value = "Sin especificar"
$ws.Range("L213").Value = "Primera"
$ws.Range("M213").Value = 760
$ws.Range("N213").Value = 6000
$ws.Range("O213").Value = 7000
$ws.Range("P213").Value = 6461
$ws.Range("Q213").Value = "$/bandeja 4 kilos"
$ws.Range("R213").Value = "Perú"
$ws.Range("S213").Value = 1615
$ws.Range("T213").Value = 4

Write-Output "Row 213 inserted and populated"
